$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Straightforward text/value updates (percent-change cells, URLs, coin names,
# and price cells that are not at risk of being auto-converted to a number).
$textUpdates = @(
    @('D2', '36.957.06'),
    @('E2', '  -0.51%  '),
    @('D3', '2.038.23'),
    @('E3', '  -0.63%  '),
    @('E4', '  -0.12%  '),
    @('E5', '  -1.79%  '),
    @('E6', '  -1.59%  '),
    @('E7', '  -1.76%  '),
    @('E8', '  -0.07%  '),
    @('E9', '  -2.03%  '),
    @('E10', '  -2.00%  '),
    @('E11', '  +2.35%  '),
    @('E12', '  -5.29%  '),
    @('E13', '  +7.85%  '),
    @('D14', '2.331.60'),
    @('E14', '  -0.91%  '),
    @('E15', '  +0.53%  '),
    @('D16', '2.026.35'),
    @('E16', '  -1.28%  '),
    @('E17', '  +4.11%  '),
    @('D18', '36.898.72'),
    @('E18', '  -0.60%  '),
    @('E19', '  -1.84%  '),
    @('D20', '0.0₃0881'),
    @('E20', '  -2.32%  '),
    @('E21', '  -0.54%  '),
    @('E22', '  -0.93%  '),
    @('E23', '  +0.05%  '),
    @('E24', '  +1.85%  '),
    @('E25', '  +2.22%  '),
    @('B26', 'PancakeSwap'),
    @('C26', 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'),
    @('E26', '  -2.10%  '),
    @('B27', 'Monero'),
    @('C27', 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'),
    @('E27', '  -0.55%  '),
    @('E28', '  -0.40%  '),
    @('E29', '  +17.83%  '),
    @('E30', '  -0.95%  '),
    @('E31', '  -4.08%  '),
    @('E32', '  +5.00%  '),
    @('E33', '  -0.98%  '),
    @('E34', '  -0.11%  '),
    @('E35', '  -4.34%  '),
    @('E37', '  -0.23%  '),
    @('E38', '  -4.85%  '),
    @('B39', 'THORChain'),
    @('C39', 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'),
    @('E39', '  -2.11%  '),
    @('B40', 'HuobiToken'),
    @('C40', 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'),
    @('E40', '  -1.85%  '),
    @('E41', '  -0.14%  '),
    @('E42', '  +0.63%  '),
    @('E43', '  -13.43%  '),
    @('E44', '  +0.81%  '),
    @('E45', '  -4.13%  '),
    @('D46', '1.291.16'),
    @('E46', '  +0.49%  '),
    @('E47', '  -4.70%  '),
    @('E48', '  -0.62%  '),
    @('B49', 'FraxShare'),
    @('C49', 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'),
    @('E49', '  -2.13%  '),
    @('B50', 'FTXToken'),
    @('C50', 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'),
    @('E50', '  +5.43%  '),
    @('B51', 'RocketPoolETH'),
    @('C51', 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'),
    @('D51', '2.216.22'),
    @('E51', '  -1.36%  '),
)
foreach ($item in $textUpdates) {
    $ws.Range($item[0]).Value = $item[1]
}

# Price cells whose new text (e.g. "245.29") parses as a plain number -- force
# the cell to Text format first so Excel keeps the literal string instead of
# silently converting it to a floating point number, then drop the format
# override again so the cell keeps its original (default) style.
$numericLookingUpdates = @(
    @('D5', '245.29'),
    @('D6', '0.654'),
    @('D7', '58.30'),
    @('D9', '0.375'),
    @('D10', '0.0768'),
    @('D12', '15.29'),
    @('D13', '0.879'),
    @('D15', '5.62'),
    @('D17', '18.19'),
    @('D19', '73.46'),
    @('D21', '5.35'),
    @('D22', '235.22'),
    @('D24', '2.44'),
    @('D25', '9.56'),
    @('D26', '2.15'),
    @('D27', '167.58'),
    @('D28', '19.86'),
    @('D29', '5.59'),
    @('D30', '0.124'),
    @('D32', '4.73'),
    @('D33', '0.0610'),
    @('D35', '0.0859'),
    @('D37', '2.24'),
    @('D38', '1.30'),
    @('D39', '5.20'),
    @('D40', '3.10'),
    @('D41', '0.0221'),
    @('D42', '1.14'),
    @('D43', '0.0949'),
    @('D44', '96.98'),
    @('D45', '16.88'),
    @('D47', '2.35'),
    @('D49', '6.66'),
    @('D50', '3.61'),
)
foreach ($item in $numericLookingUpdates) {
    $cell = $ws.Range($item[0])
    $cell.NumberFormat = "@"
    $cell.Value = $item[1]
    $cell.ClearFormats()
}
